$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: write a numeric-looking value as genuine TEXT (shared string), matching
# how the source workbook stores these "medium"/"hard" result cells, without
# disturbing the cell's existing style (border/fill/numFmt stay untouched).
function Set-TextValue($addr, $text) {
    $r = $ws.Range($addr)
    $r.Formula = '=TEXT(' + $text + ',"@")'
    $r.Copy()
    $r.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
}

# --- Row 17 (TrpB) "medium" results: S17:X17 ---
Set-TextValue "S17" '"0.01"'
Set-TextValue "T17" '"0.01"'
Set-TextValue "U17" '"2.0"'
Set-TextValue "V17" '"1.0"'
$ws.Range("W17").Value = 0
$ws.Range("X17").Value = 80

# --- Row 17 (TrpB) "hard" results: Z17:AE17 ---
Set-TextValue "Z17" '"0.0"'
Set-TextValue "AA17" '"0.0"'
Set-TextValue "AB17" '"2.0"'
Set-TextValue "AC17" '"0.0"'
$ws.Range("AD17").Value = 2
$ws.Range("AE17").Value = 80

# --- New row 19 ---
$ws.Range("AG19").Value = "S"

# --- View state: scroll + selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 17
$ws.Range("AA19").Select()
